# Refresh the cryptocurrency Price (D) and Volume/1h-change (E) columns
# with the latest scraped snapshot (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.075.46"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "1.789.98"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'323.56"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.4312"
$ws.Range("E7").Value = "  -2.86%  "
$ws.Range("D8").Value = "'0.3623"
$ws.Range("E8").Value = "  -3.15%  "
$ws.Range("D9").Value = "'44.71"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").Value = "'0.07509"
$ws.Range("E10").Value = "  -3.94%  "
$ws.Range("D11").Value = "'1.115"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "'21.62"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").Value = "'6.148"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "'7.316"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "1.786.19"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").Value = "'92.20"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").Value = "'0.06357"
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'17.23"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").Value = "'5.971"
$ws.Range("D23").Value = "28.073.82"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").Value = "'11.40"
$ws.Range("E24").Value = "  -2.49%  "
$ws.Range("D25").Value = "'2.141"
$ws.Range("E25").Value = "  -7.81%  "
$ws.Range("D26").Value = "'158.99"
$ws.Range("E26").Value = "  +3.61%  "
$ws.Range("D27").Value = "'20.36"
$ws.Range("E27").Value = "  -2.61%  "
$ws.Range("D28").Value = "1.997.27"
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("D29").Value = "'2.170"
$ws.Range("E29").Value = "  -8.72%  "
$ws.Range("D30").Value = "'127.21"
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").Value = "'1.156"
$ws.Range("E31").Value = "  -5.55%  "
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").Value = "'0.08978"
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("D34").Value = "'3.522"
$ws.Range("E34").Value = "  -4.10%  "
$ws.Range("D35").Value = "'12.60"
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("D36").Value = "'0.02315"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").Value = "'5.091"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").Value = "'0.6443"
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("D39").Value = "'0.2113"
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("D40").Value = "'0.06053"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("D41").Value = "'1.188"
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("D42").Value = "'1.422"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").Value = "'0.9996"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "'7.836"
$ws.Range("E44").Value = "  -2.77%  "
$ws.Range("D45").Value = "'13.55"
$ws.Range("E45").Value = "  -2.26%  "
$ws.Range("D46").Value = "'0.5989"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").Value = "'3.706"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").Value = "'124.65"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("D50").Value = "'1.153"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "'0.06947"
$ws.Range("E51").Value = "  +0.49%  "
